$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.6873773333333334
$ws.Range("H2").Value = 2.062132
$ws.Range("I2").Value = 0.02660947569874856
$ws.Range("J2").Value = 0.02660947569874856
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 297.8183156666666
$ws.Range("N2").Value = 893.454947
$ws.Range("O2").Value = 0.8852156413092672
$ws.Range("P2").Value = 0.8852156413092673
$ws.Range("Q2").Value = 204.7135596407782
$ws.Range("R2").Value = 1842.422036767004
$ws.Range("S2").Value = 0.02355512409557107
$ws.Range("T2").Value = 0.02355512409557107

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.6873773333333334
$ws.Range("H3").Value = 2.062132
$ws.Range("I3").Value = 0.02660947569874856
$ws.Range("J3").Value = 0.02660947569874856
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.07234764413494278
$ws.Range("P3").Value = 0.0723476441349428
$ws.Range("Q3").Value = 16.73100098026178
$ws.Range("R3").Value = 150.579008822356
$ws.Range("S3").Value = 0.001925132878470469
$ws.Range("T3").Value = 0.001925132878470469

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.6873773333333334
$ws.Range("H4").Value = 2.062132
$ws.Range("I4").Value = 0.02660947569874856
$ws.Range("J4").Value = 0.02660947569874856
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.04243671455578994
$ws.Range("P4").Value = 0.04243671455578994
$ws.Range("Q4").Value = 9.813847034296002
$ws.Range("R4").Value = 88.324623308664
$ws.Range("S4").Value = 0.001129218724707022
$ws.Range("T4").Value = 0.001129218724707022

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.913984666666667
$ws.Range("H5").Value = 17.741954
$ws.Range("I5").Value = 0.2289398029860915
$ws.Range("J5").Value = 0.2289398029860915
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 297.8183156666666
$ws.Range("N5").Value = 893.454947
$ws.Range("O5").Value = 0.8852156413092672
$ws.Range("P5").Value = 0.8852156413092673
$ws.Range("Q5").Value = 1761.29295230516
$ws.Range("R5").Value = 15851.63657074644
$ws.Range("S5").Value = 0.2026610945215502
$ws.Range("T5").Value = 0.2026610945215503

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.913984666666667
$ws.Range("H6").Value = 17.741954
$ws.Range("I6").Value = 0.2289398029860915
$ws.Range("J6").Value = 0.2289398029860915
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.07234764413494278
$ws.Range("P6").Value = 0.0723476441349428
$ws.Range("Q6").Value = 143.9484231687202
$ws.Range("R6").Value = 1295.535808518482
$ws.Range("S6").Value = 0.01656325539476166
$ws.Range("T6").Value = 0.01656325539476166

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.913984666666667
$ws.Range("H7").Value = 17.741954
$ws.Range("I7").Value = 0.2289398029860915
$ws.Range("J7").Value = 0.2289398029860915
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.04243671455578994
$ws.Range("P7").Value = 0.04243671455578994
$ws.Range("Q7").Value = 84.43534295841201
$ws.Range("R7").Value = 759.918086625708
$ws.Range("S7").Value = 0.009715453069779547
$ws.Range("T7").Value = 0.00971545306977955

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.230689
$ws.Range("H8").Value = 57.692067
$ws.Range("I8").Value = 0.7444507213151601
$ws.Range("J8").Value = 0.7444507213151601
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 297.8183156666666
$ws.Range("N8").Value = 893.454947
$ws.Range("O8").Value = 0.8852156413092672
$ws.Range("P8").Value = 0.8852156413092673
$ws.Range("Q8").Value = 5727.251407089495
$ws.Range("R8").Value = 51545.26266380545
$ws.Range("S8").Value = 0.6589994226921461
$ws.Range("T8").Value = 0.6589994226921461

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.230689
$ws.Range("H9").Value = 57.692067
$ws.Range("I9").Value = 0.7444507213151601
$ws.Range("J9").Value = 0.7444507213151601
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.34034433333333
$ws.Range("N9").Value = 73.021033
$ws.Range("O9").Value = 0.07234764413494278
$ws.Range("P9").Value = 0.0723476441349428
$ws.Range("Q9").Value = 468.0815920272457
$ws.Range("R9").Value = 4212.734328245211
$ws.Range("S9").Value = 0.05385925586171066
$ws.Range("T9").Value = 0.05385925586171068

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.230689
$ws.Range("H10").Value = 57.692067
$ws.Range("I10").Value = 0.7444507213151601
$ws.Range("J10").Value = 0.7444507213151601
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.277234
$ws.Range("N10").Value = 42.831702
$ws.Range("O10").Value = 0.04243671455578994
$ws.Range("P10").Value = 0.04243671455578994
$ws.Range("Q10").Value = 274.561046834226
$ws.Range("R10").Value = 2471.049421508034
$ws.Range("S10").Value = 0.03159204276130337
$ws.Range("T10").Value = 0.03159204276130338

